{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// The document starts with a centered date line, followed by a single\n// 20-row x 5-column table of simple arithmetic problems (e.g. \"27+4=31\").\n// This script updates the date line and replaces the text of every table\n// cell with its new value, matching the positional old->new value mapping\n// below (one pair of old/new text per cell, in row-major order). Using\n// Range.insertText(..., \"Replace\") on the existing paragraph/cell range\n// (rather than deleting/recreating runs) keeps each run's original\n// formatting (fonts, size, alignment, etc.) untouched.\n\nconst NEW_DATE = \"2024-12-13 Friday\";\nconst OLD_DATE = \"2024-12-12 Thursday\";\n\n// Expected current value of each table cell, row-major (20 rows x 5 cols).\n// Used only as a sanity check so the script fails loudly if the table\n// layout does not match what this script expects.\nconst OLD_GRID = [\n  [\"27+4=31\", \"24+18=42\", \"27+56=83\", \"18+8=26\", \"41-36=5\"],\n  [\"61-38=23\", \"29+23=52\", \"93-34=59\", \"92-46=46\", \"85-79=6\"],\n  [\"85-49=36\", \"8+89=97\", \"94-25=69\", \"81-24=57\", \"39+25=64\"],\n  [\"36+46=82\", \"16+55=71\", \"7+84=91\", \"77+14=91\", \"50-49=1\"],\n  [\"91-16=75\", \"62-43=19\", \"81-12=69\", \"4+89=93\", \"81-66=15\"],\n  [\"52-17=35\", \"6+46=52\", \"40-27=13\", \"91-15=76\", \"37+26=63\"],\n  [\"87-49=38\", \"47+6=53\", \"40-11=29\", \"34+49=83\", \"46-19=27\"],\n  [\"2+59=61\", \"57-18=39\", \"68-49=19\", \"19+77=96\", \"92-46=46\"],\n  [\"85-28=57\", \"93-87=6\", \"36+56=92\", \"16+79=95\", \"43-38=5\"],\n  [\"32-28=4\", \"92-69=23\", \"56-48=8\", \"2+79=81\", \"53-19=34\"],\n  [\"10-4=6\", \"77-68=9\", \"19+54=73\", \"63-46=17\", \"30-13=17\"],\n  [\"64-17=47\", \"47+15=62\", \"75+6=81\", \"96-18=78\", \"6+68=74\"],\n  [\"67+7=74\", \"53+28=81\", \"71-9=62\", \"16+39=55\", \"83-25=58\"],\n  [\"19+14=33\", \"9+54=63\", \"8+54=62\", \"18+49=67\", \"5+58=63\"],\n  [\"87-39=48\", \"14+58=72\", \"91-2=89\", \"48+44=92\", \"85-26=59\"],\n  [\"7+26=33\", \"59+39=98\", \"23+19=42\", \"70-8=62\", \"49+6=55\"],\n  [\"71-22=49\", \"92-16=76\", \"76-59=17\", \"71-35=36\", \"55+16=71\"],\n  [\"82-8=74\", \"63-18=45\", \"62-45=17\", \"48+15=63\", \"64+29=93\"],\n  [\"9+4=13\", \"38+16=54\", \"32-14=18\", \"6+39=45\", \"47+18=65\"],\n  [\"71-56=15\", \"47+35=82\", \"55+38=93\", \"89+3=92\", \"40-28=12\"]\n];\n\n// New value for each table cell, row-major, aligned with OLD_GRID.\nconst NEW_GRID = [\n  [\"65-57=8\", \"14+57=71\", \"79+8=87\", \"93-4=89\", \"83-34=49\"],\n  [\"28+19=47\", \"61-32=29\", \"49+22=71\", \"42-39=3\", \"39+7=46\"],\n  [\"70-17=53\", \"66+27=93\", \"81-65=16\", \"36+47=83\", \"58-19=39\"],\n  [\"13+19=32\", \"49+48=97\", \"5+38=43\", \"37+45=82\", \"56-47=9\"],\n  [\"66+28=94\", \"45-17=28\", \"45-17=28\", \"81-26=55\", \"51-37=14\"],\n  [\"36+25=61\", \"27+19=46\", \"35+57=92\", \"54+9=63\", \"39+49=88\"],\n  [\"58+37=95\", \"23-9=14\", \"24-19=5\", \"80-45=35\", \"70-32=38\"],\n  [\"52-5=47\", \"57-38=19\", \"67-59=8\", \"80-31=49\", \"24+7=31\"],\n  [\"13+28=41\", \"52-24=28\", \"78+19=97\", \"29+58=87\", \"2+89=91\"],\n  [\"7+75=82\", \"81-27=54\", \"6+57=63\", \"66-29=37\", \"48+16=64\"],\n  [\"63-29=34\", \"33+29=62\", \"8+75=83\", \"34+47=81\", \"23-14=9\"],\n  [\"62-14=48\", \"49+23=72\", \"66-9=57\", \"9+12=21\", \"53+9=62\"],\n  [\"92-74=18\", \"6+85=91\", \"61-25=36\", \"96-68=28\", \"58+34=92\"],\n  [\"84-49=35\", \"59+8=67\", \"68-29=39\", \"43-6=37\", \"83-24=59\"],\n  [\"94-49=45\", \"80-76=4\", \"55+7=62\", \"17+24=41\", \"9+47=56\"],\n  [\"15+58=73\", \"2+89=91\", \"44+48=92\", \"51-34=17\", \"93-67=26\"],\n  [\"76-48=28\", \"38+24=62\", \"36-27=9\", \"60-4=56\", \"59+17=76\"],\n  [\"24+38=62\", \"39+47=86\", \"39+53=92\", \"68-9=59\", \"26+56=82\"],\n  [\"18+28=46\", \"32-18=14\", \"80-12=68\", \"81-54=27\", \"91-63=28\"],\n  [\"40-17=23\", \"75+16=91\", \"52-9=43\", \"28+69=97\", \"9+55=64\"]\n];\n\n// --- Update the title paragraph (the date line) ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.load(\"text\");\nawait context.sync();\n\nif (titleParagraph.text.trim() !== OLD_DATE) {\n  throw new Error(\n    `Unexpected title text: \"${titleParagraph.text}\" (expected \"${OLD_DATE}\")`\n  );\n}\ntitleParagraph.getRange(\"Whole\").insertText(NEW_DATE, \"Replace\");\n\n// --- Update every cell of the first (only) table ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nfor (let r = 0; r < rowCount; r++) {\n  const colCount = table.values[r].length;\n  for (let c = 0; c < colCount; c++) {\n    const current = table.values[r][c];\n    const expected = OLD_GRID[r][c];\n    if (current !== expected) {\n      throw new Error(\n        `Unexpected text in cell (${r}, ${c}): \"${current}\" (expected \"${expected}\")`\n      );\n    }\n    const cell = table.getCell(r, c);\n    cell.body.getRange(\"Whole\").insertText(NEW_GRID[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The document starts with a centered date line, followed by a single\n# 20-row x 5-column table of simple arithmetic problems (e.g. \"27+4=31\").\n# This script updates the date line and replaces the text of every table\n# cell with its new value, matching the positional old->new value mapping\n# below (one pair of old/new text per cell, in row-major order). Assigning\n# directly to an existing Range's .Text (rather than deleting/recreating\n# the range) keeps each run's original formatting (fonts, size, alignment,\n# etc.) untouched.\n\n$NewDate = \"2024-12-13 Friday\"\n$OldDate = \"2024-12-12 Thursday\"\n\n# Expected current value of each table cell, row-major (20 rows x 5 cols).\n# Used only as a sanity check so the script fails loudly if the table\n# layout does not match what this script expects.\n$OldGrid = @(\n    @(\"27+4=31\", \"24+18=42\", \"27+56=83\", \"18+8=26\", \"41-36=5\"),\n    @(\"61-38=23\", \"29+23=52\", \"93-34=59\", \"92-46=46\", \"85-79=6\"),\n    @(\"85-49=36\", \"8+89=97\", \"94-25=69\", \"81-24=57\", \"39+25=64\"),\n    @(\"36+46=82\", \"16+55=71\", \"7+84=91\", \"77+14=91\", \"50-49=1\"),\n    @(\"91-16=75\", \"62-43=19\", \"81-12=69\", \"4+89=93\", \"81-66=15\"),\n    @(\"52-17=35\", \"6+46=52\", \"40-27=13\", \"91-15=76\", \"37+26=63\"),\n    @(\"87-49=38\", \"47+6=53\", \"40-11=29\", \"34+49=83\", \"46-19=27\"),\n    @(\"2+59=61\", \"57-18=39\", \"68-49=19\", \"19+77=96\", \"92-46=46\"),\n    @(\"85-28=57\", \"93-87=6\", \"36+56=92\", \"16+79=95\", \"43-38=5\"),\n    @(\"32-28=4\", \"92-69=23\", \"56-48=8\", \"2+79=81\", \"53-19=34\"),\n    @(\"10-4=6\", \"77-68=9\", \"19+54=73\", \"63-46=17\", \"30-13=17\"),\n    @(\"64-17=47\", \"47+15=62\", \"75+6=81\", \"96-18=78\", \"6+68=74\"),\n    @(\"67+7=74\", \"53+28=81\", \"71-9=62\", \"16+39=55\", \"83-25=58\"),\n    @(\"19+14=33\", \"9+54=63\", \"8+54=62\", \"18+49=67\", \"5+58=63\"),\n    @(\"87-39=48\", \"14+58=72\", \"91-2=89\", \"48+44=92\", \"85-26=59\"),\n    @(\"7+26=33\", \"59+39=98\", \"23+19=42\", \"70-8=62\", \"49+6=55\"),\n    @(\"71-22=49\", \"92-16=76\", \"76-59=17\", \"71-35=36\", \"55+16=71\"),\n    @(\"82-8=74\", \"63-18=45\", \"62-45=17\", \"48+15=63\", \"64+29=93\"),\n    @(\"9+4=13\", \"38+16=54\", \"32-14=18\", \"6+39=45\", \"47+18=65\"),\n    @(\"71-56=15\", \"47+35=82\", \"55+38=93\", \"89+3=92\", \"40-28=12\")\n)\n\n# New value for each table cell, row-major, aligned with $OldGrid.\n$NewGrid = @(\n    @(\"65-57=8\", \"14+57=71\", \"79+8=87\", \"93-4=89\", \"83-34=49\"),\n    @(\"28+19=47\", \"61-32=29\", \"49+22=71\", \"42-39=3\", \"39+7=46\"),\n    @(\"70-17=53\", \"66+27=93\", \"81-65=16\", \"36+47=83\", \"58-19=39\"),\n    @(\"13+19=32\", \"49+48=97\", \"5+38=43\", \"37+45=82\", \"56-47=9\"),\n    @(\"66+28=94\", \"45-17=28\", \"45-17=28\", \"81-26=55\", \"51-37=14\"),\n    @(\"36+25=61\", \"27+19=46\", \"35+57=92\", \"54+9=63\", \"39+49=88\"),\n    @(\"58+37=95\", \"23-9=14\", \"24-19=5\", \"80-45=35\", \"70-32=38\"),\n    @(\"52-5=47\", \"57-38=19\", \"67-59=8\", \"80-31=49\", \"24+7=31\"),\n    @(\"13+28=41\", \"52-24=28\", \"78+19=97\", \"29+58=87\", \"2+89=91\"),\n    @(\"7+75=82\", \"81-27=54\", \"6+57=63\", \"66-29=37\", \"48+16=64\"),\n    @(\"63-29=34\", \"33+29=62\", \"8+75=83\", \"34+47=81\", \"23-14=9\"),\n    @(\"62-14=48\", \"49+23=72\", \"66-9=57\", \"9+12=21\", \"53+9=62\"),\n    @(\"92-74=18\", \"6+85=91\", \"61-25=36\", \"96-68=28\", \"58+34=92\"),\n    @(\"84-49=35\", \"59+8=67\", \"68-29=39\", \"43-6=37\", \"83-24=59\"),\n    @(\"94-49=45\", \"80-76=4\", \"55+7=62\", \"17+24=41\", \"9+47=56\"),\n    @(\"15+58=73\", \"2+89=91\", \"44+48=92\", \"51-34=17\", \"93-67=26\"),\n    @(\"76-48=28\", \"38+24=62\", \"36-27=9\", \"60-4=56\", \"59+17=76\"),\n    @(\"24+38=62\", \"39+47=86\", \"39+53=92\", \"68-9=59\", \"26+56=82\"),\n    @(\"18+28=46\", \"32-18=14\", \"80-12=68\", \"81-54=27\", \"91-63=28\"),\n    @(\"40-17=23\", \"75+16=91\", \"52-9=43\", \"28+69=97\", \"9+55=64\")\n)\n\n$d = $word.ActiveDocument\n\n# --- Update the title paragraph (the date line) ---\n$titleParagraph = $d.Paragraphs.Item(1)\n$titleText = $titleParagraph.Range.Text.TrimEnd([char]13, [char]10)\nif ($titleText -ne $OldDate) {\n    throw \"Unexpected title text: '$titleText' (expected '$OldDate')\"\n}\n$titleParagraph.Range.Text = $NewDate\n\n# --- Update every cell of the first (only) table ---\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        # Cell.Range.Text includes trailing cell-mark characters (CR + BEL);\n        # strip them before comparing against the expected value.\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        $expected = $OldGrid[$r - 1][$c - 1]\n        if ($current -ne $expected) {\n            throw \"Unexpected text in cell ($r, $c): '$current' (expected '$expected')\"\n        }\n        $cell.Range.Text = $NewGrid[$r - 1][$c - 1]\n    }\n}\n\nWrite-Output \"done\"\n"}
